$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet pointed at a pre-production Sura environment via a live
# hyperlink on the URL cell (B2). The new environment info is now a plain
# value (no hyperlink), so drop the hyperlink entirely.
$ws.Hyperlinks.Delete()

# Replace the environment/credential values in the data row (row 2) with
# the new "gw" (Gestion Documental) environment values.
# A1/B1/C1/D1/E1 headers stay: Ambiente | URL | Usuario | Contrasenia | NroPoliza
#
# Leading apostrophes re-assert the existing "quote prefix" text formatting
# on A2 (style index 1) so the cell keeps looking/behaving like plain text,
# matching how Excel preserves that formatting when the value is edited.
$ws.Range("A2").Value = "'ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("C2").Value = "su"
$ws.Range("D2").Value = "gw"
$ws.Range("E2").Value = 4104016708

# Select the whole data row with A2 as the active cell, matching the
# selection Excel saved after editing these values.
$ws.Range("A2:E2").Select()
